# Invoice 4: apply client discount note + corrected labor quantity, then
# let Excel recalc the subtotal/tax/total formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CLIENT CODE footnote (A31): replace placeholder text with the client's
# actual contact email.
$ws.Range("A31").Value = "charlie.charlie@mail.com"

# Labor line item quantity/amount (E18 / F18 shared formula) was corrected
# from 1 to 120 hours.
$ws.Range("E18").Value = 120

# Recalculate so SUBTOTAL (F21), TAX (F23) and TOTAL (F24) pick up the
# new labor amount.
$excel.Calculate()
